$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers
$ws.Cells.Item(1,1).Value = 'Token'
$ws.Cells.Item(1,2).Value = 'Coords'
$ws.Cells.Item(1,3).Value = 'Association Ligne'
$ws.Cells.Item(1,4).Value = 'Association Col'
$ws.Cells.Item(1,5).Value = 'Moyenne'
$ws.Cells.Item(1,6).Value = 'Temps total'
$ws.Cells.Item(1,7).Value = 'Choix Final'
$ws.Cells.Item(1,8).Value = 'Id Campagne'

# Row 2
$ws.Cells.Item(2,1).Value = 'C2IDe29cc6'
$ws.Cells.Item(2,2).Value = '1:2'', 1:3'', 2:3'', 2:4'', 3:4'', 4:4'', 3:3'', 4:3'''
$ws.Cells.Item(2,3).Value = 'produit B'', produit C'', produit C'', produit C'', produit C'', produit C'', produit C'', produit C'''
$ws.Cells.Item(2,4).Value = 'Produit \xe9quitable'', Produit \xe9quitable'', u"Pr\xe9sence d''un label", u"Pr\xe9sence d''un label", produit durable'', Produit bio'', produit durable'', Produit bio'''
$ws.Cells.Item(2,5).Value = 448.5
$ws.Cells.Item(2,6).Value = 3588
$ws.Cells.Item(2,7).Value = 'produit C'
$ws.Cells.Item(2,8).Value = 2

# Row 3
$ws.Cells.Item(3,1).Value = 'C2ID32c652'
$ws.Cells.Item(3,2).Value = '1:1'', 2:2'', 1:1'', 2:3'', 2:4'', 3:4'', 4:4'', 4:3'', 1:4'''
$ws.Cells.Item(3,3).Value = 'produit C'', produit B'', produit C'', produit C'', Produit A'', Produit A'', Produit A'', produit C'', Produit A'''
$ws.Cells.Item(3,4).Value = 'produit durable'', Produit bio'', produit durable'', Produit bio'', Produit bio'', u"Pr\xe9sence d''un label", Produit \xe9quitable'', Produit \xe9quitable'', produit durable'''
$ws.Cells.Item(3,5).Value = 871.4444444444445
$ws.Cells.Item(3,6).Value = 7843
$ws.Cells.Item(3,7).Value = 'produit C'
$ws.Cells.Item(3,8).Value = 2

# Row 4
$ws.Cells.Item(4,1).Value = 'C2ID15d2c1'
$ws.Cells.Item(4,2).Value = '2:1'', 2:2'', 2:3'', 2:4'', 1:4'', 1:3'''
$ws.Cells.Item(4,3).Value = 'produit C'', produit C'', Produit A'', produit B'', produit B'', Produit A'''
$ws.Cells.Item(4,4).Value = 'Produit bio'', Produit bio'', Produit bio'', Produit bio'', u"Pr\xe9sence d''un label", u"Pr\xe9sence d''un label"'
$ws.Cells.Item(4,5).Value = 389
$ws.Cells.Item(4,6).Value = 2334
$ws.Cells.Item(4,7).Value = 'Produit A'
$ws.Cells.Item(4,8).Value = 2

# Row 5
$ws.Cells.Item(5,1).Value = 'C2IDc00290'
$ws.Cells.Item(5,2).Value = '1:1'', 2:2'', 2:3'', 3:3'', 3:4'', 4:4'', 4:3'''
$ws.Cells.Item(5,3).Value = 'produit C'', Produit A'', produit C'', produit C'', produit B'', produit B'', produit C'''
$ws.Cells.Item(5,4).Value = 'Produit bio'', produit durable'', produit durable'', u"Pr\xe9sence d''un label", u"Pr\xe9sence d''un label", Produit \xe9quitable'', Produit \xe9quitable'''
$ws.Cells.Item(5,5).Value = 1153.428571428571
$ws.Cells.Item(5,6).Value = 8074
$ws.Cells.Item(5,7).Value = 'produit C'
$ws.Cells.Item(5,8).Value = 2

# Rows 6-8 (partial, reusing row 5 text in B/C/D, blank A/E/F)
$ws.Cells.Item(6,2).Value = '1:1'', 2:2'', 2:3'', 3:3'', 3:4'', 4:4'', 4:3'''
$ws.Cells.Item(6,3).Value = 'produit C'', Produit A'', produit C'', produit C'', produit B'', produit B'', produit C'''
$ws.Cells.Item(6,4).Value = 'Produit bio'', produit durable'', produit durable'', u"Pr\xe9sence d''un label", u"Pr\xe9sence d''un label", Produit \xe9quitable'', Produit \xe9quitable'''
$ws.Cells.Item(6,7).Value = 'produit C'
$ws.Cells.Item(6,8).Value = 2
$ws.Cells.Item(6,1).Borders.LineStyle = -4142
$ws.Cells.Item(6,5).Borders.LineStyle = -4142
$ws.Cells.Item(6,6).Borders.LineStyle = -4142
$ws.Cells.Item(7,2).Value = '1:1'', 2:2'', 2:3'', 3:3'', 3:4'', 4:4'', 4:3'''
$ws.Cells.Item(7,3).Value = 'produit C'', Produit A'', produit C'', produit C'', produit B'', produit B'', produit C'''
$ws.Cells.Item(7,4).Value = 'Produit bio'', produit durable'', produit durable'', u"Pr\xe9sence d''un label", u"Pr\xe9sence d''un label", Produit \xe9quitable'', Produit \xe9quitable'''
$ws.Cells.Item(7,7).Value = 'produit C'
$ws.Cells.Item(7,8).Value = 2
$ws.Cells.Item(7,1).Borders.LineStyle = -4142
$ws.Cells.Item(7,5).Borders.LineStyle = -4142
$ws.Cells.Item(7,6).Borders.LineStyle = -4142
$ws.Cells.Item(8,2).Value = '1:1'', 2:2'', 2:3'', 3:3'', 3:4'', 4:4'', 4:3'''
$ws.Cells.Item(8,3).Value = 'produit C'', Produit A'', produit C'', produit C'', produit B'', produit B'', produit C'''
$ws.Cells.Item(8,4).Value = 'Produit bio'', produit durable'', produit durable'', u"Pr\xe9sence d''un label", u"Pr\xe9sence d''un label", Produit \xe9quitable'', Produit \xe9quitable'''
$ws.Cells.Item(8,7).Value = 'produit C'
$ws.Cells.Item(8,8).Value = 2
$ws.Cells.Item(8,1).Borders.LineStyle = -4142
$ws.Cells.Item(8,5).Borders.LineStyle = -4142
$ws.Cells.Item(8,6).Borders.LineStyle = -4142

# Row 9 (all blank, anchor dimension)
$ws.Cells.Item(9,1).Borders.LineStyle = -4142
$ws.Cells.Item(9,2).Borders.LineStyle = -4142
$ws.Cells.Item(9,3).Borders.LineStyle = -4142
$ws.Cells.Item(9,4).Borders.LineStyle = -4142
$ws.Cells.Item(9,5).Borders.LineStyle = -4142
$ws.Cells.Item(9,6).Borders.LineStyle = -4142
$ws.Cells.Item(9,7).Borders.LineStyle = -4142
$ws.Cells.Item(9,8).Borders.LineStyle = -4142

# Header styling (bold, thin border, center/top alignment)
$hdr = $ws.Range("A1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
